$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("A2").Value = 89.68085106382979
$ws.Range("B2").Value = 87.36196319018404
$ws.Range("C2").Value = 2.31888787364575
$ws.Range("D2").Value = 0.022
$ws.Range("A3").Value = 95.80851063829788
$ws.Range("B3").Value = 94.66257668711657
$ws.Range("C3").Value = 1.145933951181306
$ws.Range("D3").Value = 0.044
$ws.Range("A4").Value = 91.02127659574468
$ws.Range("B4").Value = 88.83435582822086
$ws.Range("C4").Value = 2.186920767523816
$ws.Range("D4").Value = 0.019
$ws.Range("A5").Value = 89.02127659574468
$ws.Range("B5").Value = 85.42944785276073
$ws.Range("C5").Value = 3.591828742983949
$ws.Range("D5").Value = 0
$ws.Range("A6").Value = 83.1063829787234
$ws.Range("B6").Value = 80.30674846625767
$ws.Range("C6").Value = 2.79963451246573
$ws.Range("D6").Value = 0.025
$ws.Range("A7").Value = 86.06382978723404
$ws.Range("B7").Value = 83.61963190184049
$ws.Range("C7").Value = 2.44419788539355
$ws.Range("D7").Value = 0.017
$ws.Range("A8").Value = 85.04255319148936
$ws.Range("B8").Value = 82.20858895705521
$ws.Range("C8").Value = 2.833964234434148
$ws.Range("D8").Value = 0.019
$ws.Range("A9").Value = 92.42553191489361
$ws.Range("B9").Value = 90.06134969325153
$ws.Range("C9").Value = 2.364182221642082
$ws.Range("D9").Value = 0.004
$ws.Range("A10").Value = 77.93617021276596
$ws.Range("B10").Value = 75.30674846625767
$ws.Range("C10").Value = 2.629421746508285
$ws.Range("D10").Value = 0.05
$ws.Range("A11").Value = 87.82978723404256
$ws.Range("B11").Value = 85.2760736196319
$ws.Range("C11").Value = 2.553713614410654
$ws.Range("D11").Value = 0.005
$ws.Range("A12").Value = 96.17021276595744
$ws.Range("B12").Value = 93.12883435582822
$ws.Range("C12").Value = 3.041378410129226
$ws.Range("D12").Value = 0.001
$ws.Range("A13").Value = 88.55512572533848
$ws.Range("B13").Value = 86.01784718349136
$ws.Range("C13").Value = 2.53727854184713
$ws.Range("D13").Value = 0.003
